# household_new.xlsx edit
# - Rename the data value 'instanceName' -> 'instance_name' in the
#   "survey" sheet's linked_table row (line_text.elementName column, M7).
# - The corresponding line_text.isInstanceMetadata flag (N7) flips from
#   TRUE to FALSE, matching the fact that instance_name is no longer
#   treated as instance metadata.
# - Update the active selection on the "survey" sheet from J8 (with the
#   sheet scrolled so column G is leftmost) to E7 at the default scroll
#   position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("M7").Value = "instance_name"
$ws.Range("N7").Value = $false

$ws.Activate()
$ws.Range("E7").Select()
